{"js": "// Replace the date line and each \"AxB=\" equation cell with its new value.\n// Every old value below is unique within the document, so a literal\n// (non-wildcard) search for the whole old string unambiguously identifies\n// the single run that needs to change.\nconst replacements = [\n  [\"2025-03-12 Wednesday\", \"2025-03-13 Thursday\"],\n  [\"776\u00d72=\", \"648\u00d74=\"],\n  [\"941\u00d73=\", \"636\u00d74=\"],\n  [\"101\u00d77=\", \"359\u00d79=\"],\n  [\"201\u00d75=\", \"566\u00d74=\"],\n  [\"307\u00d72=\", \"579\u00d73=\"],\n  [\"552\u00d77=\", \"228\u00d76=\"],\n  [\"467\u00d72=\", \"858\u00d73=\"],\n  [\"614\u00d72=\", \"540\u00d76=\"],\n  [\"170\u00d75=\", \"385\u00d79=\"],\n  [\"905\u00d73=\", \"433\u00d74=\"],\n  [\"672\u00d74=\", \"142\u00d74=\"],\n  [\"493\u00d78=\", \"370\u00d75=\"],\n  [\"505\u00d78=\", \"836\u00d78=\"],\n  [\"433\u00d75=\", \"636\u00d73=\"],\n  [\"354\u00d77=\", \"577\u00d78=\"],\n  [\"911\u00d79=\", \"675\u00d76=\"],\n  [\"778\u00d74=\", \"642\u00d73=\"],\n  [\"273\u00d78=\", \"627\u00d75=\"],\n  [\"821\u00d75=\", \"630\u00d79=\"],\n  [\"418\u00d74=\", \"781\u00d72=\"],\n  [\"399\u00d74=\", \"132\u00d76=\"],\n  [\"132\u00d72=\", \"533\u00d72=\"],\n  [\"474\u00d79=\", \"933\u00d75=\"],\n  [\"428\u00d76=\", \"394\u00d78=\"],\n  [\"854\u00d74=\", \"136\u00d75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the date line and each \"AxB=\" equation cell with its new value.\n# Every old value is unique within the document, so Find/Replace on the\n# literal whole string unambiguously targets the single run to edit.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-03-12 Wednesday\", \"2025-03-13 Thursday\"),\n    @(\"776\u00d72=\", \"648\u00d74=\"),\n    @(\"941\u00d73=\", \"636\u00d74=\"),\n    @(\"101\u00d77=\", \"359\u00d79=\"),\n    @(\"201\u00d75=\", \"566\u00d74=\"),\n    @(\"307\u00d72=\", \"579\u00d73=\"),\n    @(\"552\u00d77=\", \"228\u00d76=\"),\n    @(\"467\u00d72=\", \"858\u00d73=\"),\n    @(\"614\u00d72=\", \"540\u00d76=\"),\n    @(\"170\u00d75=\", \"385\u00d79=\"),\n    @(\"905\u00d73=\", \"433\u00d74=\"),\n    @(\"672\u00d74=\", \"142\u00d74=\"),\n    @(\"493\u00d78=\", \"370\u00d75=\"),\n    @(\"505\u00d78=\", \"836\u00d78=\"),\n    @(\"433\u00d75=\", \"636\u00d73=\"),\n    @(\"354\u00d77=\", \"577\u00d78=\"),\n    @(\"911\u00d79=\", \"675\u00d76=\"),\n    @(\"778\u00d74=\", \"642\u00d73=\"),\n    @(\"273\u00d78=\", \"627\u00d75=\"),\n    @(\"821\u00d75=\", \"630\u00d79=\"),\n    @(\"418\u00d74=\", \"781\u00d72=\"),\n    @(\"399\u00d74=\", \"132\u00d76=\"),\n    @(\"132\u00d72=\", \"533\u00d72=\"),\n    @(\"474\u00d79=\", \"933\u00d75=\"),\n    @(\"428\u00d76=\", \"394\u00d78=\"),\n    @(\"854\u00d74=\", \"136\u00d75=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute(\n        $oldText,\n        $false,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        $newText,\n        2\n    ) | Out-Null\n}\n"}
